$d = $word.ActiveDocument

$d.Content.Find.Execute("39+34=73", $true, $true, $false, $false, $false, $true, 1, $false, "31-1=30", 2) | Out-Null
$d.Content.Find.Execute("68-4=64", $true, $true, $false, $false, $false, $true, 1, $false, "90+3=93", 2) | Out-Null
$d.Content.Find.Execute("42+15=57", $true, $true, $false, $false, $false, $true, 1, $false, "13-3=10", 2) | Out-Null
$d.Content.Find.Execute("51+36=87", $true, $true, $false, $false, $false, $true, 1, $false, "24+41=65", 2) | Out-Null
$d.Content.Find.Execute("62+32=94", $true, $true, $false, $false, $false, $true, 1, $false, "79-69=10", 2) | Out-Null
$d.Content.Find.Execute("37+6=43", $true, $true, $false, $false, $false, $true, 1, $false, "95-37=58", 2) | Out-Null
$d.Content.Find.Execute("84-17=67", $true, $true, $false, $false, $false, $true, 1, $false, "70-63=7", 2) | Out-Null
$d.Content.Find.Execute("56+10=66", $true, $true, $false, $false, $false, $true, 1, $false, "4+75=79", 2) | Out-Null
$d.Content.Find.Execute("36+5=41", $true, $true, $false, $false, $false, $true, 1, $false, "14+76=90", 2) | Out-Null
$d.Content.Find.Execute("34+39=73", $true, $true, $false, $false, $false, $true, 1, $false, "87+0=87", 2) | Out-Null
$d.Content.Find.Execute("84-63=21", $true, $true, $false, $false, $false, $true, 1, $false, "56-40=16", 2) | Out-Null
$d.Content.Find.Execute("67-2=65", $true, $true, $false, $false, $false, $true, 1, $false, "11+72=83", 2) | Out-Null
$d.Content.Find.Execute("87-86=1", $true, $true, $false, $false, $false, $true, 1, $false, "11+42=53", 2) | Out-Null
$d.Content.Find.Execute("55-16=39", $true, $true, $false, $false, $false, $true, 1, $false, "96-76=20", 2) | Out-Null
$d.Content.Find.Execute("1+31=32", $true, $true, $false, $false, $false, $true, 1, $false, "1+23=24", 2) | Out-Null
$d.Content.Find.Execute("22+36=58", $true, $true, $false, $false, $false, $true, 1, $false, "46+50=96", 2) | Out-Null
$d.Content.Find.Execute("23+28=51", $true, $true, $false, $false, $false, $true, 1, $false, "67+30=97", 2) | Out-Null
$d.Content.Find.Execute("59-47=12", $true, $true, $false, $false, $false, $true, 1, $false, "0+85=85", 2) | Out-Null
$d.Content.Find.Execute("21+38=59", $true, $true, $false, $false, $false, $true, 1, $false, "97-86=11", 2) | Out-Null
$d.Content.Find.Execute("58-47=11", $true, $true, $false, $false, $false, $true, 1, $false, "69-20=49", 2) | Out-Null
$d.Content.Find.Execute("97-60=37", $true, $true, $false, $false, $false, $true, 1, $false, "78+6=84", 2) | Out-Null
$d.Content.Find.Execute("31+2=33", $true, $true, $false, $false, $false, $true, 1, $false, "3+54=57", 2) | Out-Null
$d.Content.Find.Execute("97-26=71", $true, $true, $false, $false, $false, $true, 1, $false, "43+44=87", 2) | Out-Null
$d.Content.Find.Execute("45-6=39", $true, $true, $false, $false, $false, $true, 1, $false, "54-54=0", 2) | Out-Null
$d.Content.Find.Execute("76-55=21", $true, $true, $false, $false, $false, $true, 1, $false, "21+51=72", 2) | Out-Null
$d.Content.Find.Execute("45+35=80", $true, $true, $false, $false, $false, $true, 1, $false, "81-60=21", 2) | Out-Null
$d.Content.Find.Execute("43+33=76", $true, $true, $false, $false, $false, $true, 1, $false, "60+7=67", 2) | Out-Null
$d.Content.Find.Execute("68-29=39", $true, $true, $false, $false, $false, $true, 1, $false, "46-12=34", 2) | Out-Null
$d.Content.Find.Execute("8+44=52", $true, $true, $false, $false, $false, $true, 1, $false, "69+7=76", 2) | Out-Null
$d.Content.Find.Execute("5+6=11", $true, $true, $false, $false, $false, $true, 1, $false, "20+15=35", 2) | Out-Null
$d.Content.Find.Execute("77-16=61", $true, $true, $false, $false, $false, $true, 1, $false, "60-23=37", 2) | Out-Null
$d.Content.Find.Execute("90-16=74", $true, $true, $false, $false, $false, $true, 1, $false, "68+5=73", 2) | Out-Null
$d.Content.Find.Execute("1+13=14", $true, $true, $false, $false, $false, $true, 1, $false, "8+10=18", 2) | Out-Null
$d.Content.Find.Execute("17-13=4", $true, $true, $false, $false, $false, $true, 1, $false, "74-74=0", 2) | Out-Null
$d.Content.Find.Execute("40+5=45", $true, $true, $false, $false, $false, $true, 1, $false, "57+40=97", 2) | Out-Null
$d.Content.Find.Execute("58-12=46", $true, $true, $false, $false, $false, $true, 1, $false, "8+24=32", 2) | Out-Null
$d.Content.Find.Execute("80-73=7", $true, $true, $false, $false, $false, $true, 1, $false, "36-24=12", 2) | Out-Null
$d.Content.Find.Execute("62-51=11", $true, $true, $false, $false, $false, $true, 1, $false, "0+4=4", 2) | Out-Null
$d.Content.Find.Execute("24+56=80", $true, $true, $false, $false, $false, $true, 1, $false, "49+6=55", 2) | Out-Null
$d.Content.Find.Execute("76-64=12", $true, $true, $false, $false, $false, $true, 1, $false, "92-49=43", 2) | Out-Null
$d.Content.Find.Execute("24+49=73", $true, $true, $false, $false, $false, $true, 1, $false, "76-29=47", 2) | Out-Null
$d.Content.Find.Execute("57+41=98", $true, $true, $false, $false, $false, $true, 1, $false, "8+15=23", 2) | Out-Null
$d.Content.Find.Execute("53-21=32", $true, $true, $false, $false, $false, $true, 1, $false, "87-31=56", 2) | Out-Null
$d.Content.Find.Execute("80-34=46", $true, $true, $false, $false, $false, $true, 1, $false, "17-4=13", 2) | Out-Null
$d.Content.Find.Execute("39-17=22", $true, $true, $false, $false, $false, $true, 1, $false, "5+27=32", 2) | Out-Null
$d.Content.Find.Execute("82-75=7", $true, $true, $false, $false, $false, $true, 1, $false, "68-26=42", 2) | Out-Null
$d.Content.Find.Execute("11-4=7", $true, $true, $false, $false, $false, $true, 1, $false, "46+44=90", 2) | Out-Null
$d.Content.Find.Execute("89-72=17", $true, $true, $false, $false, $false, $true, 1, $false, "61-27=34", 2) | Out-Null
$d.Content.Find.Execute("76-35=41", $true, $true, $false, $false, $false, $true, 1, $false, "68-30=38", 2) | Out-Null
$d.Content.Find.Execute("77-43=34", $true, $true, $false, $false, $false, $true, 1, $false, "77-17=60", 2) | Out-Null
$d.Content.Find.Execute("83-62=21", $true, $true, $false, $false, $false, $true, 1, $false, "30+66=96", 2) | Out-Null
$d.Content.Find.Execute("19+12=31", $true, $true, $false, $false, $false, $true, 1, $false, "82-39=43", 2) | Out-Null
$d.Content.Find.Execute("1+54=55", $true, $true, $false, $false, $false, $true, 1, $false, "1+23=24", 2) | Out-Null
$d.Content.Find.Execute("36+29=65", $true, $true, $false, $false, $false, $true, 1, $false, "33-26=7", 2) | Out-Null
$d.Content.Find.Execute("36+20=56", $true, $true, $false, $false, $false, $true, 1, $false, "41+3=44", 2) | Out-Null
$d.Content.Find.Execute("94-40=54", $true, $true, $false, $false, $false, $true, 1, $false, "28+36=64", 2) | Out-Null
$d.Content.Find.Execute("73-53=20", $true, $true, $false, $false, $false, $true, 1, $false, "54-31=23", 2) | Out-Null
$d.Content.Find.Execute("73-0=73", $true, $true, $false, $false, $false, $true, 1, $false, "71-23=48", 2) | Out-Null
$d.Content.Find.Execute("74+2=76", $true, $true, $false, $false, $false, $true, 1, $false, "74-28=46", 2) | Out-Null
$d.Content.Find.Execute("36-30=6", $true, $true, $false, $false, $false, $true, 1, $false, "22+48=70", 2) | Out-Null
$d.Content.Find.Execute("14+19=33", $true, $true, $false, $false, $false, $true, 1, $false, "71-53=18", 2) | Out-Null
$d.Content.Find.Execute("25-25=0", $true, $true, $false, $false, $false, $true, 1, $false, "66-49=17", 2) | Out-Null
$d.Content.Find.Execute("19-13=6", $true, $true, $false, $false, $false, $true, 1, $false, "85-53=32", 2) | Out-Null
$d.Content.Find.Execute("54-34=20", $true, $true, $false, $false, $false, $true, 1, $false, "28-5=23", 2) | Out-Null
$d.Content.Find.Execute("98-20=78", $true, $true, $false, $false, $false, $true, 1, $false, "69-9=60", 2) | Out-Null
$d.Content.Find.Execute("65-39=26", $true, $true, $false, $false, $false, $true, 1, $false, "88+5=93", 2) | Out-Null
$d.Content.Find.Execute("47-35=12", $true, $true, $false, $false, $false, $true, 1, $false, "8+14=22", 2) | Out-Null
$d.Content.Find.Execute("85+1=86", $true, $true, $false, $false, $false, $true, 1, $false, "12+34=46", 2) | Out-Null
$d.Content.Find.Execute("82-8=74", $true, $true, $false, $false, $false, $true, 1, $false, "65+29=94", 2) | Out-Null
$d.Content.Find.Execute("54+27=81", $true, $true, $false, $false, $false, $true, 1, $false, "25-16=9", 2) | Out-Null
$d.Content.Find.Execute("97-81=16", $true, $true, $false, $false, $false, $true, 1, $false, "33-15=18", 2) | Out-Null
$d.Content.Find.Execute("91-40=51", $true, $true, $false, $false, $false, $true, 1, $false, "28-15=13", 2) | Out-Null
$d.Content.Find.Execute("23+4=27", $true, $true, $false, $false, $false, $true, 1, $false, "24+10=34", 2) | Out-Null
$d.Content.Find.Execute("99-72=27", $true, $true, $false, $false, $false, $true, 1, $false, "53+0=53", 2) | Out-Null
$d.Content.Find.Execute("81-11=70", $true, $true, $false, $false, $false, $true, 1, $false, "82-49=33", 2) | Out-Null
$d.Content.Find.Execute("20+73=93", $true, $true, $false, $false, $false, $true, 1, $false, "46+17=63", 2) | Out-Null
$d.Content.Find.Execute("77-56=21", $true, $true, $false, $false, $false, $true, 1, $false, "25+41=66", 2) | Out-Null
$d.Content.Find.Execute("70+29=99", $true, $true, $false, $false, $false, $true, 1, $false, "97-71=26", 2) | Out-Null
$d.Content.Find.Execute("15+84=99", $true, $true, $false, $false, $false, $true, 1, $false, "63+23=86", 2) | Out-Null
$d.Content.Find.Execute("7+78=85", $true, $true, $false, $false, $false, $true, 1, $false, "13+77=90", 2) | Out-Null
$d.Content.Find.Execute("73-38=35", $true, $true, $false, $false, $false, $true, 1, $false, "42-15=27", 2) | Out-Null
$d.Content.Find.Execute("37-11=26", $true, $true, $false, $false, $false, $true, 1, $false, "36+28=64", 2) | Out-Null
$d.Content.Find.Execute("67-55=12", $true, $true, $false, $false, $false, $true, 1, $false, "14+14=28", 2) | Out-Null
$d.Content.Find.Execute("55-6=49", $true, $true, $false, $false, $false, $true, 1, $false, "42-3=39", 2) | Out-Null
$d.Content.Find.Execute("7-7=0", $true, $true, $false, $false, $false, $true, 1, $false, "12+27=39", 2) | Out-Null
$d.Content.Find.Execute("37+50=87", $true, $true, $false, $false, $false, $true, 1, $false, "13+79=92", 2) | Out-Null
$d.Content.Find.Execute("62-48=14", $true, $true, $false, $false, $false, $true, 1, $false, "46+38=84", 2) | Out-Null
$d.Content.Find.Execute("4+2=6", $true, $true, $false, $false, $false, $true, 1, $false, "34+46=80", 2) | Out-Null
$d.Content.Find.Execute("29+67=96", $true, $true, $false, $false, $false, $true, 1, $false, "37+34=71", 2) | Out-Null
$d.Content.Find.Execute("87-2=85", $true, $true, $false, $false, $false, $true, 1, $false, "95-76=19", 2) | Out-Null
$d.Content.Find.Execute("6+47=53", $true, $true, $false, $false, $false, $true, 1, $false, "88-32=56", 2) | Out-Null
$d.Content.Find.Execute("47-45=2", $true, $true, $false, $false, $false, $true, 1, $false, "40+56=96", 2) | Out-Null
$d.Content.Find.Execute("41-2=39", $true, $true, $false, $false, $false, $true, 1, $false, "60+26=86", 2) | Out-Null
$d.Content.Find.Execute("55-42=13", $true, $true, $false, $false, $false, $true, 1, $false, "54-22=32", 2) | Out-Null
$d.Content.Find.Execute("87-72=15", $true, $true, $false, $false, $false, $true, 1, $false, "63-60=3", 2) | Out-Null
$d.Content.Find.Execute("3+11=14", $true, $true, $false, $false, $false, $true, 1, $false, "83-35=48", 2) | Out-Null
$d.Content.Find.Execute("97-21=76", $true, $true, $false, $false, $false, $true, 1, $false, "49-16=33", 2) | Out-Null
$d.Content.Find.Execute("50-8=42", $true, $true, $false, $false, $false, $true, 1, $false, "28+57=85", 2) | Out-Null
$d.Content.Find.Execute("41+43=84", $true, $true, $false, $false, $false, $true, 1, $false, "61-57=4", 2) | Out-Null
$d.Content.Find.Execute("36-25=11", $true, $true, $false, $false, $false, $true, 1, $false, "50-11=39", 2) | Out-Null
